# Week 13 logging: add new player row (N.McCrary) to the RB stats sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("RB")

# New row right after the last existing player (row 8 -> row 9)
$row = 9

$ws.Cells.Item($row, 1).Value = "N.McCrary"
for ($col = 2; $col -le 10; $col++) {
    $ws.Cells.Item($row, $col).Value = 0
}

# Move the active selection, matching the post-edit workbook state.
$ws.Range("J10").Select()
